# Final commit before merge
# Applies:
#  - remove the "Straight Arrow Connector 6" shape from the sheet's drawing
#  - clear the validation note text out of D11 (drops the now-unused shared string)
#  - remove the custom pink fill from the D11:E11 style (back to "No Fill")
#  - move the active selection from B9:C9 to F12
#  - reset the page setup back to the sheet defaults

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the straight arrow connector shape that pointed at the D column notes.
$shp = $ws.Shapes.Item("Straight Arrow Connector 6")
$shp.Delete()

# Clear the "only 10 digits..." note text from D11 (E11 stays blank, already empty).
$ws.Range("D11").Value = ""

# Drop the custom fill color on D11:E11, going back to "No Fill".
$ws.Range("D11:E11").Interior.ColorIndex = -4142
$ws.Range("D11:E11").Interior.Pattern = -4142

# Update the saved selection/active cell for the sheet view.
$ws.Range("F12").Select() | Out-Null

# Reset page setup back to the worksheet's base defaults.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
